# Fix the "Source" label typo for 2014 onshore rows: "onshore" -> "Onshore"
# (to match the capitalization used for other years), and leave the
# selection on C6 (the first corrected cell), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6:C10").Value = "Onshore"

$ws.Range("C6").Select()
